$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'249.28"
$ws.Range("D3").Value = "'21.66"
$ws.Range("D4").Value = "'5.575"
$ws.Range("D5").Value = "'0.05670"
$ws.Range("D6").Value = "'3.378"
$ws.Range("D7").Value = "'6.444"
$ws.Range("D8").Value = "'0.8055"
$ws.Range("D9").Value = "'1.041"
$ws.Range("D10").Value = "'0.1427"
$ws.Range("D11").Value = "'0.07267"
$ws.Range("D13").Value = "'0.02921"
$ws.Range("D14").Value = "'0.09273"
$ws.Range("D15").Value = "'0.001670"
$ws.Range("D16").Value = "'3.216"
$ws.Range("D18").Value = "'0.0005815"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006459"
$ws.Range("D21").Value = "'0.001050"
$ws.Range("D23").Value = "'3.976"
$ws.Range("D24").Value = "'2.113"
$ws.Range("D25").Value = "'0.3293"
$ws.Range("D27").Value = "'0.0003203"
$ws.Range("D40").Value = "'0.04128"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006916"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.003503"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1044"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.008520"
$ws.Range("D45").Value = "'0.00005643"
$ws.Range("D47").Value = "'0.7859"
$ws.Range("D48").Value = "'0.01676"
